$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price table rows per the latest data refresh.
# Column D (Price) cells are plain text in the source data (e.g. "30.255.91",
# "1.0000"), so force text format before assigning to avoid Excel silently
# re-interpreting them as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.255.91"
$ws.Range("E2").Value = "  -0.68%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.856.52"
$ws.Range("E3").Value = "  -2.29%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.16"
$ws.Range("E5").Value = "  -0.99%  "

$ws.Range("E7").Value = "  -2.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2799"
$ws.Range("E8").Value = "  -4.08%  "

$ws.Range("E9").Value = "  -3.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.858.44"
$ws.Range("E10").Value = "  -2.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07379"
$ws.Range("E11").Value = "  +0.77%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.23"
$ws.Range("E12").Value = "  -4.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.016"
$ws.Range("E13").Value = "  -3.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.05"
$ws.Range("E14").Value = "  -1.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6438"
$ws.Range("E15").Value = "  -3.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.209.30"
$ws.Range("E16").Value = "  -0.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.0000"
$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("E18").Value = "  -2.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007554"
$ws.Range("E19").Value = "  -3.92%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "223.02"
$ws.Range("E20").Value = "  +14.28%  "

$ws.Range("B21").Value = "BinanceUSD"
$ws.Range("C21").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.265"
$ws.Range("E22").Value = "  -3.33%  "

$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.053"
$ws.Range("E23").Value = "  -1.21%  "

$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.174"
$ws.Range("E24").Value = "  -3.20%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "163.38"
$ws.Range("E25").Value = "  +0.11%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.49"
$ws.Range("E26").Value = "  +0.80%  "

$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.925"
$ws.Range("E27").Value = "  -0.72%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.434"
$ws.Range("E28").Value = "  -3.20%  "

$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09197"
$ws.Range("E29").Value = "  +0.57%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.232"
$ws.Range("E30").Value = "  -2.04%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.962"
$ws.Range("E31").Value = "  -3.68%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04948"
$ws.Range("E32").Value = "  -4.12%  "

$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.141"
$ws.Range("E33").Value = "  +3.26%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7235"
$ws.Range("E34").Value = "  -1.98%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.687"
$ws.Range("E35").Value = "  -1.55%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01828"
$ws.Range("E36").Value = "  -0.64%  "

$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.595"
$ws.Range("E37").Value = "  -2.99%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8984"
$ws.Range("E38").Value = "  -2.72%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.035"
$ws.Range("E39").Value = "  -1.53%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.897"
$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "105.99"
$ws.Range("E41").Value = "  -0.84%  "

$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").Value = "  +0.54%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4238"
$ws.Range("E43").Value = "  -3.71%  "

$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.253"
$ws.Range("E44").Value = "  -4.06%  "

$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1294"
$ws.Range("E45").Value = "  -5.45%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.22"
$ws.Range("E46").Value = "  -8.12%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.484"
$ws.Range("E47").Value = "  +6.02%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.747"
$ws.Range("E48").Value = "  -3.08%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.70"
$ws.Range("E49").Value = "  -3.44%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05635"
$ws.Range("E50").Value = "  -3.39%  "

$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3782"
$ws.Range("E51").Value = "  -3.61%  "
